$wb = $excel.ActiveWorkbook

# The "Options" sheet holds the lookup lists that back the dropdown data
# validations on the "Template" sheet. This change revises the accepted
# soil-depth buckets in column K: existing values are re-labelled to
# shallower 20cm bands, and a new "1m+" option is appended in K6. The
# Template sheet's data validation range for column K is extended from
# K1:K5 to K1:K6 to include the new option.

$optionsWs = $wb.Worksheets.Item("Options")
$templateWs = $wb.Worksheets.Item("Template")

# Update existing soil depth labels in column K of the Options sheet.
$optionsWs.Range("K1").Value = "0-20cm"
$optionsWs.Range("K2").Value = "20cm-40cm"
$optionsWs.Range("K3").Value = "40cm-60cm"
$optionsWs.Range("K4").Value = "60cm-80cm"
$optionsWs.Range("K5").Value = "80cm-1m"
# New option added in the previously-empty K6 cell.
$optionsWs.Range("K6").Value = "1m+"

# Extend the data validation list range for column K on the Template
# sheet so it covers the new Options!K6 entry (was K1:K5).
$kRange = $templateWs.Range("K2:K101")
$kRange.Validation.Modify(3, 1, 1, "=Options!`$K`$1:`$K`$6")
